# "Updated slides for second training"
#  - Slide-master date stamps: 19.06.2024 -> 06.07.2024 (appears on both
#    master designs that carry that shape)
#  - Slide 1 subtitle: "18.06.2024, Daniel Krämer" -> "09.07.2024, Daniel
#    Krämer" (only the date portion is retyped)

$p = $ppt.ActivePresentation

$oldDate1 = "19.06.2024"
$newDate1 = "06.07.2024"

# The dated "Rectangle 6" shape lives on more than one slide master
# (design); walk every design's master and fix it wherever it shows up.
for ($i = 1; $i -le $p.Designs.Count; $i++) {
    $design = $p.Designs.Item($i)
    $master = $design.SlideMaster
    for ($j = 1; $j -le $master.Shapes.Count; $j++) {
        $shp = $master.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate1) {
                $tr.Text = $newDate1
            }
        }
    }
}

# Slide 1: retype just the leading date run, leaving ", Daniel Krämer" in
# its own (second) run.
$oldLead = "18.06.2024"
$newLead = "09.07.2024"
$slide1 = $p.Slides.Item(1)
for ($k = 1; $k -le $slide1.Shapes.Count; $k++) {
    $shp2 = $slide1.Shapes.Item($k)
    if ($shp2.HasTextFrame) {
        $full = $shp2.TextFrame.TextRange.Text
        if ($full.Length -ge $oldLead.Length -and $full.Substring(0, $oldLead.Length) -eq $oldLead) {
            $datePart = $shp2.TextFrame.TextRange.Characters(1, $oldLead.Length)
            $datePart.Text = $newLead
        }
    }
}
